# Fill in the previously-blank cells of the Topics/Topics compatibility
# matrix on Sheet1 with "yes" / "cant use" values (and one "cant usw" typo
# cell), then restore the active selection as it was left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("R8").Value = "cant use"
$ws.Range("S8").Value = "yes"

$ws.Range("L9:Q9").Value = "yes"
$ws.Range("R9").Value = "cant use"
$ws.Range("S9:T9").Value = "yes"

$ws.Range("J10:Q10").Value = "yes"
$ws.Range("R10").Value = "cant use"
$ws.Range("S10:T10").Value = "yes"

$ws.Range("J11:Q11").Value = "yes"
$ws.Range("R11").Value = "cant use"
$ws.Range("S11:T11").Value = "yes"

$ws.Range("I12:Q12").Value = "yes"
$ws.Range("R12").Value = "cant use"
$ws.Range("S12:T12").Value = "yes"

$ws.Range("I13:Q13").Value = "yes"
$ws.Range("R13").Value = "cant use"
$ws.Range("S13:T13").Value = "yes"

$ws.Range("I14:Q14").Value = "yes"
$ws.Range("R14").Value = "cant use"
$ws.Range("S14:T14").Value = "yes"

$ws.Range("I15:Q15").Value = "yes"
$ws.Range("R15").Value = "cant use"
$ws.Range("S15:T15").Value = "yes"

$ws.Range("I16:T16").Value = "yes"

$ws.Range("I17:T17").Value = "yes"

$ws.Range("E18").Value = "yes"
$ws.Range("I18:L18").Value = "cant use"
$ws.Range("M18").Value = "cant usw"
$ws.Range("N18:O18").Value = "cant use"
$ws.Range("P18:T18").Value = "yes"

$ws.Range("I19:T19").Value = "yes"

$ws.Range("I20:T20").Value = "yes"

# Restore the selection left by the author when the file was saved.
$ws.Range("N29").Select()
